# Add simple Buttons for GUI
# Update the "Newest Case" form-number counters on row 3 of Sheet 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("B3").Value = 1024
$ws.Range("C3").Value = 2002
$ws.Range("D3").Value = 3002
$ws.Range("E3").Value = 4002
